# Generate Report for Handoff
# Adds a new handed-off file ("eaa5d274-36f7-4cf2-81b0-0c290cbbfbe6.md") as a
# new row to each of the three worksheets (Overview, zh-cn, de-de), mirroring
# the existing row for "5a7cce6c-4826-4621-af79-f55d4e859c7b.md".

$wb = $excel.ActiveWorkbook

# Blue color used by the workbook's existing "HyperLink" style (RGB 6495ED),
# expressed as a VBA-style BGR packed value for the Font.Color property.
$hyperlinkColor = 15570276

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "eaa5d274-36f7-4cf2-81b0-0c290cbbfbe6.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-16 04:38:46"

$hOverview = $wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d86297633b48963ff92a40df831b1f646afc11f5/e2e/eaa5d274-36f7-4cf2-81b0-0c290cbbfbe6.md",
    "",
    "",
    "e2e\eaa5d274-36f7-4cf2-81b0-0c290cbbfbe6.md"
)
$wsOverview.Range("B3").Font.Underline = 2
$wsOverview.Range("B3").Font.Color = $hyperlinkColor

$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "eaa5d274-36f7-4cf2-81b0-0c290cbbfbe6.66670ad6c3e6a29c3d14727df0231c13469fe0bd.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-16 04:38:42"
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("O3").Value = "False"

$hZhCn = $wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d86297633b48963ff92a40df831b1f646afc11f5/e2e/eaa5d274-36f7-4cf2-81b0-0c290cbbfbe6.md",
    "",
    "",
    "eaa5d274-36f7-4cf2-81b0-0c290cbbfbe6.md"
)
$wsZhCn.Range("A3").Font.Underline = 2
$wsZhCn.Range("A3").Font.Color = $hyperlinkColor

$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "eaa5d274-36f7-4cf2-81b0-0c290cbbfbe6.66670ad6c3e6a29c3d14727df0231c13469fe0bd.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-16 04:38:46"
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("O3").Value = "False"

$hDeDe = $wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d86297633b48963ff92a40df831b1f646afc11f5/e2e/eaa5d274-36f7-4cf2-81b0-0c290cbbfbe6.md",
    "",
    "",
    "eaa5d274-36f7-4cf2-81b0-0c290cbbfbe6.md"
)
$wsDeDe.Range("A3").Font.Underline = 2
$wsDeDe.Range("A3").Font.Color = $hyperlinkColor

$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

Write-Output "Handback row added for eaa5d274-36f7-4cf2-81b0-0c290cbbfbe6.md on Overview, zh-cn, de-de."
